# Add a new "Images" column to the AddProduct sheet, inserted right
# before the existing "ImportDate" column (old column L), pushing
# ImportDate from L to M.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column at position 12 (L). This shifts the former
# column L ("ImportDate", with all of its per-row formatting) to
# column M, matching the diff's row/column shift.
$ws.Columns.Item(12).Insert()

# Header text for the newly inserted column.
$ws.Range("L1").Value = "Images"

# Give the new column the same width as column K (Classify), like in
# the target workbook where columns K and L share one width.
# (Inserting the column already copies K's per-cell formatting onto
# the new column automatically, so no extra style assignment needed.)
$kWidth = $ws.Columns.Item(11).ColumnWidth()
$ws.Columns.Item(12).ColumnWidth = $kWidth

# Update the view so the new column is visible / selected, matching
# the author's recorded selection in the workbook.
$ws.Range("L2").Select()
$excel.ActiveWindow.ScrollColumn = 11
